$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.241.30"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "2.523.63"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.535"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.557"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.73%  "
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "2.921.05"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").Value = "2.524.66"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "48.075.27"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "0.0₃0949"
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.84%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0794"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "120.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("D45").Value = "2.021.43"
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.50%  "
$ws.Range("E47").Value = "  +6.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.47%  "
